$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Maximum amount" (column D) values keyed by currency (column C):
#   ДЕН (denar) -> 2500, ЕУР (euro) -> 500, УСД (usd) -> 300
# Rows 2-13 already contain an (empty but styled) D cell, so just fill in the value.
$ws.Range("D2").Value  = 2500
$ws.Range("D3").Value  = 500
$ws.Range("D4").Value  = 300
$ws.Range("D5").Value  = 2500
$ws.Range("D6").Value  = 500
$ws.Range("D7").Value  = 300
$ws.Range("D8").Value  = 2500
$ws.Range("D9").Value  = 500
$ws.Range("D10").Value = 300
$ws.Range("D11").Value = 2500
$ws.Range("D12").Value = 500
$ws.Range("D13").Value = 300

# Rows 14-23 had no D cell at all before; give them the "Normal 2" cell style
# (the same style family used elsewhere in the column) and the value.
$ws.Range("D14").Style = "Normal 2"
$ws.Range("D14").Value = 2500
$ws.Range("D15").Style = "Normal 2"
$ws.Range("D15").Value = 500
$ws.Range("D16").Style = "Normal 2"
$ws.Range("D16").Value = 2500
$ws.Range("D17").Style = "Normal 2"
$ws.Range("D17").Value = 500
$ws.Range("D18").Style = "Normal 2"
$ws.Range("D18").Value = 2500
$ws.Range("D19").Style = "Normal 2"
$ws.Range("D19").Value = 500
$ws.Range("D20").Style = "Normal 2"
$ws.Range("D20").Value = 300
$ws.Range("D21").Style = "Normal 2"
$ws.Range("D21").Value = 2500
$ws.Range("D22").Style = "Normal 2"
$ws.Range("D22").Value = 500
$ws.Range("D23").Style = "Normal 2"
$ws.Range("D23").Value = 300

# The author's last selection on the sheet moved from C24 to D24.
[void]$ws.Range("D24").Select()
